$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.136.27"
$ws.Range("E2").Value = "  +2.97%  "
$ws.Range("D3").Value = "2.642.17"
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'595.76"
$ws.Range("D6").Value = "'156.31"
$ws.Range("E6").Value = "  +4.16%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +1.17%  "
$ws.Range("D9").Value = "'0.118"
$ws.Range("E9").Value = "  +7.67%  "
$ws.Range("E10").Value = "  +4.33%  "
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("E12").Value = "  +2.01%  "
$ws.Range("D13").Value = "'29.00"
$ws.Range("E13").Value = "  +5.51%  "
$ws.Range("D14").Value = "3.113.47"
$ws.Range("E14").Value = "  +2.28%  "
$ws.Range("D15").Value = "'0.0000183"
$ws.Range("E15").Value = "  +18.39%  "
$ws.Range("D16").Value = "65.032.40"
$ws.Range("E16").Value = "  +3.11%  "
$ws.Range("D17").Value = "2.617.29"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").Value = "'12.57"
$ws.Range("E18").Value = "  +2.64%  "
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("D20").Value = "'353.92"
$ws.Range("E20").Value = "  +2.66%  "
$ws.Range("D21").Value = "'7.29"
$ws.Range("E21").Value = "  +6.23%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").Value = "'68.22"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("D25").Value = "'9.54"
$ws.Range("E25").Value = "  +3.03%  "
$ws.Range("D26").Value = "'1.64"
$ws.Range("E26").Value = "  -1.58%  "
$ws.Range("D27").Value = "'8.16"
$ws.Range("E27").Value = "  +1.36%  "
$ws.Range("D28").Value = "'0.164"
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "0.0₃0946"
$ws.Range("E30").Value = "  +11.96%  "
$ws.Range("D31").Value = "'520.41"
$ws.Range("E31").Value = "  -7.93%  "
$ws.Range("E32").Value = "  +3.38%  "
$ws.Range("D33").Value = "'1.78"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("D34").Value = "'5.64"
$ws.Range("E34").Value = "  +8.14%  "
$ws.Range("D35").Value = "'6.31"
$ws.Range("E35").Value = "  +3.37%  "
$ws.Range("E36").Value = "  +3.80%  "
$ws.Range("D37").Value = "'164.75"
$ws.Range("E37").Value = "  -1.18%  "
$ws.Range("E38").Value = "  +5.72%  "
$ws.Range("D39").Value = "'20.25"
$ws.Range("E39").Value = "  +4.14%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "'42.24"
$ws.Range("E42").Value = "  +6.91%  "
$ws.Range("D43").Value = "'165.17"
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("D44").Value = "'4.10"
$ws.Range("E44").Value = "  +3.78%  "
$ws.Range("D46").Value = "'22.97"
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("E47").Value = "  +6.97%  "
$ws.Range("D48").Value = "'0.648"
$ws.Range("E48").Value = "  +3.28%  "
$ws.Range("E49").Value = "  +1.23%  "
$ws.Range("D50").Value = "'0.0985"
$ws.Range("E50").Value = "  +2.47%  "
$ws.Range("D51").Value = "'19.47"
$ws.Range("E51").Value = "  +2.35%  "
